$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp label in A1
$ws.Range("A1").Value = 'Datos actualizados a 28 de Junio de 2020 a las 02:30'

# Country data rows 4-219: country name (A) and stats (B:H)
# Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$rows = @(
    @{Row=4; A='Estados Unidos'; B=2596403; C=43447; D=1080132; E=1388119; F=0; G=512; H=128152},
    @{Row=5; A='Brasil'; B=1315941; C=35887; D=715905; E=542933; F=0; G=994; H=57103},
    @{Row=6; A='Rusia'; B=627646; C=6852; D=393352; E=225325; F=0; G=188; H=8969},
    @{Row=7; A='India'; B=529577; C=20131; D=310146; E=203328; F=0; G=414; H=16103},
    @{Row=8; A='Reino Unido'; B=310250; C=890; D=0; E=0; F=0; G=100; H=43514},
    @{Row=9; A='España'; B=295549; C=564; D=0; E=0; F=0; G=3; H=28341},
    @{Row=10; A='Peru'; B=275989; C=3625; D=164024; E=102830; F=0; G=196; H=9135},
    @{Row=11; A='Chile'; B=267766; C=4406; D=228055; E=34364; F=0; G=279; H=5347},
    @{Row=12; A='Italia'; B=240136; C=175; D=188584; E=16836; F=0; G=8; H=34716},
    @{Row=13; A='Iran'; B=220180; C=2456; D=180661; E=29155; F=0; G=125; H=10364},
    @{Row=14; A='Mexico'; B=208392; C=5441; D=120562; E=62051; F=0; G=719; H=25779},
    @{Row=15; A='Pakistan'; B=198883; C=3138; D=86906; E=107942; F=0; G=73; H=4035},
    @{Row=16; A='Turquia'; B=195883; C=1372; D=169182; E=21619; F=0; G=17; H=5082},
    @{Row=17; A='Alemania'; B=194689; C=290; D=177500; E=8163; F=0; G=0; H=9026},
    @{Row=18; A='Arabia Saudita'; B=178504; C=3927; D=122128; E=54865; F=0; G=37; H=1511},
    @{Row=19; A='Francia'; B=162936; C=0; D=75649; E=57509; F=0; G=0; H=29778},
    @{Row=20; A='Banglades'; B=133978; C=3504; D=54318; E=77965; F=0; G=34; H=1695},
    @{Row=21; A='Sudafrica'; B=131800; C=7210; D=67094; E=62293; F=0; G=73; H=2413},
    @{Row=22; A='Canada'; B=103032; C=238; D=65973; E=28543; F=0; G=8; H=8516},
    @{Row=23; A='Catar'; B=93663; C=879; D=77225; E=16328; F=0; G=1; H=110},
    @{Row=24; A='Colombia'; B=88591; C=4149; D=36273; E=49379; F=0; G=128; H=2939},
    @{Row=25; A='China'; B=83483; C=21; D=78444; E=405; F=0; G=0; H=4634},
    @{Row=26; A='Suecia'; B=65137; C=0; D=0; E=0; F=0; G=0; H=5280},
    @{Row=27; A='Egipto'; B=63923; C=1168; D=17140; E=44075; F=0; G=88; H=2708},
    @{Row=28; A='Belgica'; B=61209; C=103; D=16941; E=34536; F=0; G=1; H=9732},
    @{Row=29; A='Bielorrusia'; B=61095; C=382; D=44126; E=16592; F=0; G=4; H=377},
    @{Row=30; A='Argentina'; B=57744; C=2401; D=19143; E=37394; F=0; G=23; H=1207},
    @{Row=31; A='Ecuador'; B=54574; C=718; D=26920; E=23230; F=0; G=18; H=4424},
    @{Row=32; A='Indonesia'; B=52812; C=1385; D=21909; E=28183; F=0; G=37; H=2720},
    @{Row=33; A='Paises Bajos'; B=50074; C=69; D=0; E=0; F=0; G=2; H=6105},
    @{Row=34; A='Emiratos Arabes Unidos'; B=47360; C=387; D=35834; E=11215; F=0; G=1; H=311},
    @{Row=35; A='Kuwait'; B=44391; C=688; D=34586; E=9461; F=0; G=3; H=344},
    @{Row=36; A='Irak'; B=43262; C=2069; D=19938; E=21664; F=0; G=101; H=1660},
    @{Row=37; A='Singapur'; B=43246; C=291; D=37163; E=6057; F=0; G=0; H=26},
    @{Row=38; A='Ucrania'; B=42065; C=948; D=18701; E=22254; F=0; G=24; H=1110},
    @{Row=39; A='Portugal'; B=41189; C=323; D=26864; E=12764; F=0; G=6; H=1561},
    @{Row=40; A='Oman'; B=36953; C=919; D=20363; E=16431; F=0; G=6; H=159},
    @{Row=41; A='Filipinas'; B=34803; C=730; D=9430; E=24137; F=0; G=12; H=1236},
    @{Row=42; A='Polonia'; B=33714; C=319; D=19972; E=12307; F=0; G=6; H=1435},
    @{Row=43; A='Suiza'; B=31555; C=69; D=29100; E=493; F=0; G=0; H=1962},
    @{Row=44; A='Panama'; B=30658; C=753; D=15370; E=14696; F=0; G=17; H=592},
    @{Row=45; A='Republica Dominicana'; B=30619; C=855; D=16666; E=13235; F=0; G=6; H=718},
    @{Row=46; A='Afganistan'; B=30616; C=165; D=10674; E=19239; F=0; G=20; H=703},
    @{Row=47; A='Bolivia'; B=29423; C=920; D=7736; E=20753; F=0; G=21; H=934},
    @{Row=48; A='Rumania'; B=26022; C=325; D=18530; E=5903; F=0; G=10; H=1589},
    @{Row=49; A='Irlanda'; B=25437; C=23; D=23364; E=339; F=0; G=4; H=1734},
    @{Row=50; A='Barein'; B=25267; C=462; D=19781; E=5408; F=0; G=5; H=78},
    @{Row=51; A='Nigeria'; B=24077; C=779; D=8625; E=14894; F=0; G=4; H=558},
    @{Row=52; A='Armenia'; B=23909; C=662; D=12911; E=10577; F=0; G=11; H=421},
    @{Row=53; A='Israel'; B=23421; C=621; D=17002; E=6102; F=0; G=3; H=317},
    @{Row=54; A='Kazajistan'; B=20319; C=569; D=12738; E=7415; F=0; G=16; H=166},
    @{Row=55; A='Japon'; B=18297; C=100; D=16452; E=874; F=0; G=2; H=971},
    @{Row=56; A='Austria'; B=17580; C=58; D=16371; E=509; F=0; G=2; H=700},
    @{Row=57; A='Ghana'; B=16431; C=597; D=12257; E=4071; F=0; G=0; H=103},
    @{Row=58; A='Moldavia'; B=16080; C=304; D=8963; E=6596; F=0; G=6; H=521},
    @{Row=59; A='Honduras'; B=15994; C=628; D=1678; E=13845; F=0; G=45; H=471},
    @{Row=60; A='Azerbaiyan'; B=15890; C=521; D=8719; E=6978; F=0; G=6; H=193},
    @{Row=61; A='Guatemala'; B=15828; C=209; D=3028; E=12128; F=0; G=49; H=672},
    @{Row=62; A='Serbia'; B=13792; C=227; D=12338; E=1187; F=0; G=2; H=267},
    @{Row=63; A='Argelia'; B=12968; C=283; D=9202; E=2874; F=0; G=7; H=892},
    @{Row=64; A='Dinamarca'; B=12675; C=0; D=11508; E=563; F=0; G=0; H=604},
    @{Row=65; A='Corea del Sur'; B=12653; C=51; D=11317; E=1054; F=0; G=0; H=282},
    @{Row=66; A='Camerun'; B=12592; C=0; D=10100; E=2179; F=0; G=0; H=313},
    @{Row=67; A='Nepal'; B=12309; C=554; D=2834; E=9447; F=0; G=1; H=28},
    @{Row=68; A='Marruecos'; B=11877; C=244; D=8723; E=2934; F=0; G=2; H=220},
    @{Row=69; A='Chequia'; B=11298; C=260; D=7682; E=3267; F=0; G=0; H=349},
    @{Row=70; A='Sudan'; B=9257; C=0; D=4014; E=4671; F=0; G=0; H=572},
    @{Row=71; A='Costa de Marfil'; B=8944; C=205; D=3722; E=5156; F=0; G=2; H=66},
    @{Row=72; A='Noruega'; B=8846; C=14; D=8138; E=459; F=0; G=0; H=249},
    @{Row=73; A='Malasia'; B=8616; C=10; D=8308; E=187; F=0; G=0; H=121},
    @{Row=74; A='Uzbekistan'; B=7682; C=255; D=5240; E=2422; F=0; G=0; H=20},
    @{Row=75; A='Australia'; B=7641; C=46; D=6979; E=558; F=0; G=0; H=104},
    @{Row=76; A='Finlandia'; B=7198; C=7; D=6600; E=270; F=0; G=0; H=328},
    @{Row=77; A='Consejo Danes para los Refugiados'; B=6690; C=138; D=937; E=5600; F=0; G=4; H=153},
    @{Row=78; A='Senegal'; B=6459; C=105; D=4255; E=2102; F=0; G=4; H=102},
    @{Row=79; A='Republica de Macedonia'; B=5906; C=148; D=2236; E=3393; F=0; G=9; H=277},
    @{Row=80; A='Kenia'; B=5811; C=278; D=1936; E=3734; F=0; G=4; H=141},
    @{Row=81; A='Tayikistan'; B=5799; C=52; D=4391; E=1356; F=0; G=0; H=52},
    @{Row=82; A='El Salvador'; B=5727; C=210; D=3447; E=2137; F=0; G=10; H=143},
    @{Row=83; A='Haiti'; B=5722; C=179; D=641; E=4983; F=0; G=2; H=98},
    @{Row=84; A='Etiopia'; B=5570; C=145; D=2015; E=3461; F=0; G=5; H=94},
    @{Row=85; A='Guinea'; B=5291; C=31; D=4215; E=1046; F=0; G=1; H=30},
    @{Row=86; A='Gabon'; B=5209; C=0; D=2327; E=2842; F=0; G=0; H=40},
    @{Row=87; A='Venezuela'; B=4779; C=0; D=1327; E=3411; F=0; G=0; H=41},
    @{Row=88; A='Republica de Yibuti'; B=4643; C=0; D=4348; E=243; F=0; G=0; H=52},
    @{Row=89; A='Bulgaria'; B=4625; C=112; D=2475; E=1934; F=0; G=1; H=216},
    @{Row=90; A='Kirguistan'; B=4513; C=309; D=2212; E=2255; F=0; G=3; H=46},
    @{Row=91; A='Luxemburgo'; B=4217; C=44; D=3978; E=129; F=0; G=0; H=110},
    @{Row=92; A='Hungria'; B=4138; C=11; D=2681; E=879; F=0; G=0; H=578},
    @{Row=93; A='Mauritania'; B=4025; C=118; D=1344; E=2560; F=0; G=1; H=121},
    @{Row=94; A='Bosnia y Herzegovina'; B=3935; C=0; D=2338; E=1419; F=0; G=0; H=178},
    @{Row=95; A='Guayana Francesa'; B=3461; C=191; D=1249; E=2200; F=0; G=1; H=12},
    @{Row=96; A='Republica de Africa Central'; B=3429; C=89; D=699; E=2685; F=0; G=5; H=45},
    @{Row=97; A='Grecia'; B=3366; C=23; D=1374; E=1801; F=0; G=0; H=191},
    @{Row=98; A='Tailandia'; B=3162; C=0; D=3053; E=51; F=0; G=0; H=58},
    @{Row=99; A='Costa Rica'; B=2979; C=143; D=1325; E=1641; F=0; G=1; H=13},
    @{Row=100; A='Somalia'; B=2878; C=0; D=868; E=1920; F=0; G=0; H=90},
    @{Row=101; A='Croacia'; B=2624; C=85; D=2152; E=365; F=0; G=0; H=107},
    @{Row=102; A='Mayotte'; B=2508; C=0; D=2218; E=258; F=0; G=0; H=32},
    @{Row=103; A='Albania'; B=2330; C=61; D=1346; E=931; F=0; G=2; H=53},
    @{Row=104; A='Cuba'; B=2330; C=5; D=2187; E=57; F=0; G=1; H=86},
    @{Row=105; A='Maldivas'; B=2305; C=22; D=1875; E=422; F=0; G=0; H=8},
    @{Row=106; A='Nicaragua'; B=2170; C=0; D=1238; E=858; F=0; G=0; H=74},
    @{Row=107; A='Mali'; B=2118; C=58; D=1398; E=607; F=0; G=0; H=113},
    @{Row=108; A='Sri Lanka'; B=2033; C=19; D=1639; E=383; F=0; G=0; H=11},
    @{Row=109; A='Madagascar'; B=2005; C=83; D=907; E=1082; F=0; G=0; H=16},
    @{Row=110; A='Guinea Ecuatorial'; B=2001; C=0; D=515; E=1454; F=0; G=0; H=32},
    @{Row=111; A='Estonia'; B=1986; C=0; D=1812; E=105; F=0; G=0; H=69},
    @{Row=112; A='Sudan del Sur'; B=1942; C=0; D=224; E=1682; F=0; G=0; H=36},
    @{Row=113; A='Paraguay'; B=1942; C=231; D=1045; E=882; F=0; G=2; H=15},
    @{Row=114; A='Islandia'; B=1836; C=4; D=1814; E=12; F=0; G=0; H=10},
    @{Row=115; A='Estado de Palestina'; B=1815; C=258; D=446; E=1365; F=0; G=1; H=4},
    @{Row=116; A='Lituania'; B=1813; C=5; D=1503; E=232; F=0; G=0; H=78},
    @{Row=117; A='Libano'; B=1719; C=22; D=1153; E=533; F=0; G=0; H=33},
    @{Row=118; A='Eslovaquia'; B=1657; C=14; D=1455; E=174; F=0; G=0; H=28},
    @{Row=119; A='Guinea-Bisau'; B=1614; C=0; D=317; E=1275; F=0; G=0; H=22},
    @{Row=120; A='Eslovenia'; B=1572; C=14; D=1376; E=87; F=0; G=0; H=109},
    @{Row=121; A='Zambia'; B=1531; C=0; D=1233; E=277; F=0; G=0; H=21},
    @{Row=122; A='Nueva Zelanda'; B=1522; C=2; D=1484; E=16; F=0; G=0; H=22},
    @{Row=123; A='Sierra Leona'; B=1410; C=16; D=937; E=414; F=0; G=0; H=59},
    @{Row=124; A='Hong Kong'; B=1198; C=1; D=1095; E=96; F=0; G=0; H=7},
    @{Row=125; A='Tunez'; B=1168; C=4; D=1025; E=93; F=0; G=0; H=50},
    @{Row=126; A='Benin'; B=1124; C=71; D=295; E=815; F=0; G=0; H=14},
    @{Row=127; A='Letonia'; B=1115; C=3; D=932; E=153; F=0; G=0; H=30},
    @{Row=128; A='Jordania'; B=1111; C=7; D=841; E=261; F=0; G=0; H=9},
    @{Row=129; A='Yemen'; B=1103; C=14; D=417; E=390; F=0; G=3; H=296},
    @{Row=130; A='Cabo Verde'; B=1091; C=64; D=568; E=511; F=0; G=2; H=12},
    @{Row=131; A='Congo'; B=1087; C=0; D=456; E=594; F=0; G=0; H=37},
    @{Row=132; A='Niger'; B=1062; C=3; D=924; E=71; F=0; G=0; H=67},
    @{Row=133; A='Malaui'; B=1038; C=33; D=260; E=765; F=0; G=0; H=13},
    @{Row=134; A='Republica de Chipre'; B=994; C=2; D=824; E=151; F=0; G=0; H=19},
    @{Row=135; A='Burkina Faso'; B=941; C=0; D=830; E=58; F=0; G=0; H=53},
    @{Row=136; A='Uruguay'; B=924; C=5; D=818; E=80; F=0; G=0; H=26},
    @{Row=137; A='Georgia'; B=921; C=2; D=781; E=126; F=0; G=0; H=14},
    @{Row=138; A='Ruanda'; B=878; C=20; D=413; E=463; F=0; G=0; H=2},
    @{Row=139; A='Republica del Chad'; B=865; C=0; D=778; E=13; F=0; G=0; H=74},
    @{Row=140; A='Principado de Andorra'; B=855; C=0; D=799; E=4; F=0; G=0; H=52},
    @{Row=141; A='Uganda'; B=848; C=15; D=761; E=87; F=0; G=0; H=0},
    @{Row=142; A='Mozambique'; B=816; C=0; D=223; E=588; F=0; G=0; H=5},
    @{Row=143; A='Suazilandia'; B=745; C=17; D=370; E=367; F=0; G=0; H=8},
    @{Row=144; A='Liberia'; B=729; C=45; D=291; E=404; F=0; G=0; H=34},
    @{Row=145; A='Libia'; B=727; C=14; D=171; E=538; F=0; G=0; H=18},
    @{Row=146; A='Santo Tome y Principe'; B=713; C=1; D=219; E=481; F=0; G=0; H=13},
    @{Row=147; A='Crucero'; B=712; C=0; D=651; E=48; F=0; G=0; H=13},
    @{Row=148; A='San Marino'; B=698; C=0; D=656; E=0; F=0; G=0; H=42},
    @{Row=149; A='Jamaica'; B=686; C=2; D=539; E=137; F=0; G=0; H=10},
    @{Row=150; A='Malta'; B=670; C=0; D=635; E=26; F=0; G=0; H=9},
    @{Row=151; A='Togo'; B=615; C=24; D=396; E=205; F=0; G=0; H=14},
    @{Row=152; A='Zimbabue'; B=567; C=6; D=142; E=419; F=0; G=0; H=6},
    @{Row=153; A='Reunion'; B=520; C=3; D=472; E=46; F=0; G=0; H=2},
    @{Row=154; A='Tanzania'; B=509; C=0; D=183; E=305; F=0; G=0; H=21},
    @{Row=155; A='Montenegro'; B=469; C=30; D=315; E=145; F=0; G=0; H=9},
    @{Row=156; A='Taiwan'; B=447; C=0; D=435; E=5; F=0; G=0; H=7},
    @{Row=157; A='Surinam'; B=433; C=44; D=191; E=231; F=0; G=1; H=11},
    @{Row=158; A='Vietnam'; B=355; C=2; D=330; E=25; F=0; G=0; H=0},
    @{Row=159; A='Mauricio'; B=341; C=0; D=326; E=5; F=0; G=0; H=10},
    @{Row=160; A='Isla de Man'; B=336; C=0; D=312; E=0; F=0; G=0; H=24},
    @{Row=161; A='Birmania'; B=293; C=0; D=216; E=71; F=0; G=0; H=6},
    @{Row=162; A='Comoras'; B=272; C=0; D=161; E=104; F=0; G=0; H=7},
    @{Row=163; A='Angola'; B=259; C=47; D=81; E=168; F=0; G=0; H=10},
    @{Row=164; A='Siria'; B=256; C=1; D=102; E=145; F=0; G=1; H=9},
    @{Row=165; A='Martinica'; B=242; C=0; D=98; E=130; F=0; G=0; H=14},
    @{Row=166; A='Guyana'; B=230; C=15; D=109; E=109; F=0; G=0; H=12},
    @{Row=167; A='Mongolia'; B=219; C=0; D=175; E=44; F=0; G=0; H=0},
    @{Row=168; A='Islas Caimanes'; B=196; C=0; D=186; E=9; F=0; G=0; H=1},
    @{Row=169; A='Eritrea'; B=191; C=24; D=53; E=138; F=0; G=0; H=0},
    @{Row=170; A='Islas Feroe'; B=187; C=0; D=187; E=0; F=0; G=0; H=0},
    @{Row=171; A='Guadalupe'; B=182; C=0; D=157; E=11; F=0; G=0; H=14},
    @{Row=172; A='Gibraltar'; B=177; C=1; D=176; E=1; F=0; G=0; H=0},
    @{Row=173; A='Burundi'; B=170; C=26; D=115; E=54; F=0; G=0; H=1},
    @{Row=174; A='Bermudas'; B=146; C=0; D=134; E=3; F=0; G=0; H=9},
    @{Row=175; A='Brunei'; B=141; C=0; D=138; E=0; F=0; G=0; H=3},
    @{Row=176; A='Camboya'; B=139; C=9; D=129; E=10; F=0; G=0; H=0},
    @{Row=177; A='Namibia'; B=136; C=15; D=22; E=114; F=0; G=0; H=0},
    @{Row=178; A='Trinidad yTobago'; B=126; C=2; D=109; E=9; F=0; G=0; H=8},
    @{Row=179; A='Bahamas'; B=104; C=0; D=87; E=6; F=0; G=0; H=11},
    @{Row=180; A='Monaco'; B=103; C=1; D=95; E=4; F=0; G=0; H=4},
    @{Row=181; A='Aruba'; B=101; C=0; D=98; E=0; F=0; G=0; H=3},
    @{Row=182; A='Barbados'; B=97; C=0; D=90; E=0; F=0; G=0; H=7},
    @{Row=183; A='Botsuana'; B=92; C=0; D=25; E=66; F=0; G=0; H=1},
    @{Row=184; A='Liechtenstein'; B=82; C=0; D=81; E=0; F=0; G=0; H=1},
    @{Row=185; A='San Martin (Parte Holandesa)'; B=77; C=0; D=62; E=0; F=0; G=0; H=15},
    @{Row=186; A='Butan'; B=75; C=5; D=38; E=37; F=0; G=0; H=0},
    @{Row=187; A='Antigua y Barbuda'; B=65; C=0; D=22; E=40; F=0; G=0; H=3},
    @{Row=188; A='Polinesia Francesa'; B=60; C=0; D=60; E=0; F=0; G=0; H=0},
    @{Row=189; A='Macao'; B=46; C=0; D=45; E=1; F=0; G=0; H=0},
    @{Row=190; A='Gambia'; B=44; C=1; D=26; E=16; F=0; G=0; H=2},
    @{Row=191; A='San Martin (Parte Francesa)'; B=42; C=0; D=36; E=3; F=0; G=0; H=3},
    @{Row=192; A='Puerto Rico'; B=39; C=0; D=1; E=36; F=0; G=0; H=2},
    @{Row=193; A='Guam'; B=32; C=0; D=0; E=31; F=0; G=0; H=1},
    @{Row=194; A='San Vicente y las Granadinas'; B=29; C=0; D=29; E=0; F=0; G=0; H=0},
    @{Row=195; A='Lesoto'; B=24; C=0; D=4; E=20; F=0; G=0; H=0},
    @{Row=196; A='Belice'; B=24; C=1; D=18; E=4; F=0; G=0; H=2},
    @{Row=197; A='Timor Oriental'; B=24; C=0; D=24; E=0; F=0; G=0; H=0},
    @{Row=198; A='Curazao'; B=23; C=0; D=19; E=3; F=0; G=0; H=1},
    @{Row=199; A='Granada'; B=23; C=0; D=23; E=0; F=0; G=0; H=0},
    @{Row=200; A='Nueva Caledonia'; B=21; C=0; D=21; E=0; F=0; G=0; H=0},
    @{Row=201; A='Laos'; B=19; C=0; D=19; E=0; F=0; G=0; H=0},
    @{Row=202; A='Santa Lucia'; B=19; C=0; D=19; E=0; F=0; G=0; H=0},
    @{Row=203; A='Fiyi'; B=18; C=0; D=18; E=0; F=0; G=0; H=0},
    @{Row=204; A='Dominica'; B=18; C=0; D=18; E=0; F=0; G=0; H=0},
    @{Row=205; A='Islas Virgenes de los Estados Unidos'; B=17; C=0; D=0; E=17; F=0; G=0; H=0},
    @{Row=206; A='Islas Turcas y Caicos'; B=16; C=0; D=11; E=4; F=0; G=0; H=1},
    @{Row=207; A='San Cristobal y Nieves'; B=15; C=0; D=15; E=0; F=0; G=0; H=0},
    @{Row=208; A='Groenlandia'; B=13; C=0; D=13; E=0; F=0; G=0; H=0},
    @{Row=209; A='Islas Malvinas'; B=13; C=0; D=13; E=0; F=0; G=0; H=0},
    @{Row=210; A='Santa Sede'; B=12; C=0; D=12; E=0; F=0; G=0; H=0},
    @{Row=211; A='Papua Nueva Guinea'; B=11; C=0; D=8; E=3; F=0; G=0; H=0},
    @{Row=212; A='Seychelles'; B=11; C=0; D=11; E=0; F=0; G=0; H=0},
    @{Row=213; A='Montserrat'; B=11; C=0; D=10; E=0; F=0; G=0; H=1},
    @{Row=214; A='Sahara Occidental'; B=10; C=0; D=8; E=1; F=0; G=0; H=1},
    @{Row=215; A='Islas Virgenes Britanicas'; B=8; C=0; D=7; E=0; F=0; G=0; H=1},
    @{Row=216; A='Bonaire, San Eustaquio y Saba'; B=7; C=0; D=7; E=0; F=0; G=0; H=0},
    @{Row=217; A='San Bartolome'; B=6; C=0; D=6; E=0; F=0; G=0; H=0},
    @{Row=218; A='Anguila'; B=3; C=0; D=3; E=0; F=0; G=0; H=0},
    @{Row=219; A='San Pedro y Miquelon'; B=1; C=0; D=1; E=0; F=0; G=0; H=0}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
}